$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to append (distance, MEAN, STD, MIN, MAX, COUNT, Month)
$newRows = @(
    @(1100, 20.79930686950684, 2.995251417160034, 3.875397920608521, 29.5693244934082, 18476, "03"),
    @(1200, 20.72659873962402, 3.061853170394897, 3.844643115997314, 29.18489074707031, 18568, "03"),
    @(1300, 20.66297912597656, 3.103049278259277, 3.417493343353271, 29.30961799621582, 18438, "03"),
    @(1400, 20.60930633544922, 3.153537034988403, 3.420910358428955, 30.0135612487793, 18366, "03"),
    @(1500, 20.6399040222168, 3.150049686431885, 3.520009279251099, 30.31769180297852, 18392, "03")
)

$startRow = 21
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    $ws.Cells.Item($r, 5).Value = $rowData[4]
    $ws.Cells.Item($r, 6).Value = $rowData[5]
    $gCell = $ws.Cells.Item($r, 7)
    $gCell.NumberFormat = "@"
    $gCell.Value = $rowData[6]
}
